$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, pushing the existing rows 84-102 down to 85-103.
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record.
$ws.Range("A84").Value = 11
$ws.Range("B84").Value = "Vega Monumental Concepción"
$ws.Range("C84").Value = "Bíobío"
$ws.Range("D84").Value = 44754
$ws.Range("E84").Value = 8
$ws.Range("F84").Value = 100112024
$ws.Range("G84").Value = "Choclo"
$ws.Range("H84").Value = "Dulce o Americano"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 100
$ws.Range("K84").Value = 37000
$ws.Range("L84").Value = 38000
$ws.Range("M84").Value = 37500
$ws.Range("N84").Value = "$/malla 70 unidades"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 536
$ws.Range("Q84").Value = 70
$ws.Range("R84").Value = "Hortaliza"
